$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (pushes current rows 10-26 down to 11-27)
$ws.Rows("10:10").Insert()

$newDate = (Get-Date -Year 2023 -Month 9 -Day 26).Date

$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = $newDate
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = 100112039
$ws.Cells.Item(10, 7).Value = "Ciboulette"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 100
$ws.Cells.Item(10, 11).Value = 2500
$ws.Cells.Item(10, 12).Value = 2500
$ws.Cells.Item(10, 13).Value = 2500
$ws.Cells.Item(10, 14).Value = "`$/docena de atados"
$ws.Cells.Item(10, 15).Value = "Región Metropolitana"
$ws.Cells.Item(10, 16).Value = 833
$ws.Cells.Item(10, 17).Value = 3
$ws.Cells.Item(10, 18).Value = "Hortaliza"
